$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 23:47"

# Refresh per-country COVID figures (Casos totales, Nuevos casos, Casos
# activos, Recuperados, Casos criticos, Muertes hoy, Muertes). Some rows
# also change which country they show, because the underlying data is
# ranked by case count and a handful of countries swapped rank.
$ws.Range("B4").Value = 3279602
$ws.Range("C4").Value = 59603
$ws.Range("D4").Value = 1450355
$ws.Range("E4").Value = 1692771
$ws.Range("G4").Value = 654
$ws.Range("H4").Value = 136476
$ws.Range("B5").Value = 1800827
$ws.Range("C5").Value = 41724
$ws.Range("E5").Value = 544833
$ws.Range("G5").Value = 1144
$ws.Range("H5").Value = 70398
$ws.Range("B8").Value = 319646
$ws.Range("C8").Value = 3198
$ws.Range("D8").Value = 210638
$ws.Range("E8").Value = 97508
$ws.Range("G8").Value = 186
$ws.Range("H8").Value = 11500
$ws.Range("B19").Value = 199588
$ws.Range("C19").Value = 390
$ws.Range("E19").Value = 6458
$ws.Range("G19").Value = 5
$ws.Range("H19").Value = 9130
$ws.Range("B23").Value = 107023
$ws.Range("C23").Value = 218
$ws.Range("D23").Value = 70819
$ws.Range("E23").Value = 27445
$ws.Range("B28").Value = 74898
$ws.Range("C28").Value = 163
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = 5526
$ws.Range("B50").Value = 32039
$ws.Range("C50").Value = 511
$ws.Range("D50").Value = 27213
$ws.Range("E50").Value = 4722
$ws.Range("B71").Value = 12027
$ws.Range("C71").Value = 463
$ws.Range("D71").Value = 7530
$ws.Range("E71").Value = 4443
$ws.Range("A85").Value = "Bulgaria"
$ws.Range("B85").Value = 6964
$ws.Range("C85").Value = 292
$ws.Range("D85").Value = 3308
$ws.Range("E85").Value = 3389
$ws.Range("G85").Value = 5
$ws.Range("H85").Value = 267
$ws.Range("A86").Value = "Costa Rica"
$ws.Range("B86").Value = 6845
$ws.Range("C86").Value = 360
$ws.Range("D86").Value = 2110
$ws.Range("E86").Value = 4709
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 26
$ws.Range("B91").Value = 5942
$ws.Range("C91").Value = 71
$ws.Range("D91").Value = 3004
$ws.Range("E91").Value = 2892
$ws.Range("A107").Value = "Paraguay"
$ws.Range("B107").Value = 2736
$ws.Range("C107").Value = 98
$ws.Range("D107").Value = 1256
$ws.Range("E107").Value = 1460
$ws.Range("H107").Value = 20
$ws.Range("A108").Value = "Mayotte"
$ws.Range("B108").Value = 2711
$ws.Range("C108").Value = 9
$ws.Range("D108").Value = 2480
$ws.Range("E108").Value = 194
$ws.Range("H108").Value = 37
$ws.Range("B137").Value = 1099
$ws.Range("C137").Value = 2
$ws.Range("D137").Value = 978
$ws.Range("A138").Value = "Burkina Faso"
$ws.Range("B138").Value = 1020
$ws.Range("C138").Value = 15
$ws.Range("D138").Value = 862
$ws.Range("E138").Value = 105
$ws.Range("H138").Value = 53
$ws.Range("A139").Value = "Montenegro"
$ws.Range("B139").Value = 1019
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 320
$ws.Range("E139").Value = 680
$ws.Range("A140").Value = "Republica de Chipre"
$ws.Range("B140").Value = 1013
$ws.Range("C140").Value = 3
$ws.Range("D140").Value = 839
$ws.Range("E140").Value = 155
$ws.Range("H140").Value = 19
$ws.Range("A141").Value = "Uganda"
$ws.Range("B141").Value = 1006
$ws.Range("C141").Value = 6
$ws.Range("D141").Value = 938
$ws.Range("E141").Value = 68
$ws.Range("H141").Value = 0
$ws.Range("B168").Value = 255
$ws.Range("C168").Value = 6
$ws.Range("E168").Value = 142
$ws.Range("G168").Value = 1
$ws.Range("H168").Value = 15
$ws.Range("A173").Value = "Guadalupe"
$ws.Range("B173").Value = 190
$ws.Range("C173").Value = 6
$ws.Range("D173").Value = 157
$ws.Range("E173").Value = 19
$ws.Range("H173").Value = 14
$ws.Range("A174").Value = "Islas Feroe"
$ws.Range("B174").Value = 188
$ws.Range("D174").Value = 188
$ws.Range("E174").Value = 0
$ws.Range("H174").Value = 0
